$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update week number (column B) for existing rows 2-13 from 38 -> 41
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 2).Value = 41
}

# Update task text (column C) for rows 7-13 to the new/reassigned task strings
$ws.Cells.Item(7, 3).Value  = "Prepare data untuk Update Bundling Complete Me & Hyangyu"
$ws.Cells.Item(8, 3).Value  = "Setting Shoppefood MOMOYO"
$ws.Cells.Item(9, 3).Value  = "Instalasi Kamera baru & Peremajaan CCTV Toko Mas An An"
$ws.Cells.Item(10, 3).Value = "Reposisi Kamera CCTV Complete Me"
$ws.Cells.Item(11, 3).Value = "Follow up update penambahan PPN Ke Helpdesk Nagatech"
$ws.Cells.Item(12, 3).Value = "Trial Update Siatem NSI Toko Mas An An & Complete Mulia"
$ws.Cells.Item(13, 3).Value = "Update Sistem NSI Toko Mas an an & Complete Mulia"

# Add two new rows (14-15) following the same pattern as rows 7-13
$ws.Cells.Item(14, 1).Value = 2023
$ws.Cells.Item(14, 2).Value = 41
$ws.Cells.Item(14, 3).Value = "Setting Perubahan Menu & Harga Baru di MOKA Pos Complete Me & Hyangyu"
$ws.Cells.Item(14, 3).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "NON"

$ws.Cells.Item(15, 1).Value = 2023
$ws.Cells.Item(15, 2).Value = 41
$ws.Cells.Item(15, 3).Value = "Revisi Bundling Menu di ID Jurnal Complete Me & Hyangyu"
$ws.Cells.Item(15, 3).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "NON"

# Update the selection to reflect the new active cell (cosmetic)
$ws.Range("C12").Select()
